$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 994, pushing the existing rows 994-1017 down to
# 1000-1023 (matches the dimension growing from A1:T1017 to A1:T1023).
$ws.Range("A994:A999").EntireRow.Insert()

# Populate the 6 freshly-inserted rows (994-999) with the new weekly price
# records. Columns A,B,C,E,F,G,H,I,J,K are constant for this market/product
# subset (Vega Central Mapocho de Santiago - Limón, "Sin especificar").

$newRows = @(
    @{ Row = 994; D = 44448; L = "1a amarillo"; M = 460; N = 3700; O = 4000; P = 3843; Q = "`$/malla 18 kilos"; R = "Provincia de Melipilla";    S = 214; T = 18 },
    @{ Row = 995; D = 44448; L = "1a amarillo"; M = 430; N = 3700; O = 4000; P = 3874; Q = "`$/malla 18 kilos"; R = "Región de O'Higgins";       S = 215; T = 18 },
    @{ Row = 996; D = 44448; L = "2a amarillo"; M = 540; N = 2700; O = 3000; P = 2833; Q = "`$/malla 18 kilos"; R = "Provincia de Melipilla";    S = 157; T = 18 },
    @{ Row = 997; D = 44448; L = "2a amarillo"; M = 420; N = 2700; O = 3000; P = 2871; Q = "`$/malla 18 kilos"; R = "Región de O'Higgins";       S = 160; T = 18 },
    @{ Row = 998; D = 44448; L = "3a amarillo"; M = 500; N = 1700; O = 2000; P = 1868; Q = "`$/malla 18 kilos"; R = "Provincia de Melipilla";    S = 104; T = 18 },
    @{ Row = 999; D = 44448; L = "3a amarillo"; M = 540; N = 1800; O = 2000; P = 1911; Q = "`$/malla 18 kilos"; R = "Región de O'Higgins";       S = 106; T = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
